# Apply the edits described by the diff to Sheet1.
# The commit adds a new location ("Indianapolis EMS Conference Center"),
# removes "INDIANA HOSA", shifting subsequent rows up, updates a few
# "user_ratings_total" counts, and refreshes the raw index column (A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 27
# Row 3
$ws.Range("A3").Value = 14
# Row 4
$ws.Range("A4").Value = 52
# Row 10
$ws.Range("A10").Value = 29
$ws.Range("E10").Value = 473
# Row 12
$ws.Range("A12").Value = 1
# Row 14
$ws.Range("A14").Value = 44
# Row 17
$ws.Range("A17").Value = 16
# Row 20
$ws.Range("A20").Value = 38
$ws.Range("C20").Value = "Indiana Catholic Conference"
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
# Row 21
$ws.Range("A21").Value = 19
$ws.Range("C21").Value = "Indiana Convention Center"
$ws.Range("D21").Value = 4.5
$ws.Range("E21").Value = 528
# Row 22
$ws.Range("A22").Value = 46
$ws.Range("C22").Value = "Indiana Government Center Conference Rooms A and B"
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
# Row 23
$ws.Range("A23").Value = 25
$ws.Range("C23").Value = "Indiana United Methodist Conference"
$ws.Range("D23").Value = 4.5
$ws.Range("E23").Value = 4
# Row 24
$ws.Range("A24").Value = 58
$ws.Range("C24").Value = "Indiana Water Environment Association"
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
# Row 25
$ws.Range("A25").Value = 5
$ws.Range("C25").Value = "Indiana Wesleyan University - Greenwood Education and Conference Center"
$ws.Range("D25").Value = 4.3
$ws.Range("E25").Value = 7
# Row 26
$ws.Range("A26").Value = 7
$ws.Range("C26").Value = "Indiana Wesleyan University - Indianapolis North Education and Conference Center"
$ws.Range("D26").Value = 4.8
$ws.Range("E26").Value = 12
# Row 27
$ws.Range("A27").Value = 6
$ws.Range("C27").Value = "Indiana Wesleyan University - Indianapolis West Education and Conference Center"
$ws.Range("D27").Value = 4.3
$ws.Range("E27").Value = 3
# Row 28
$ws.Range("A28").Value = 34
$ws.Range("C28").Value = "Indiana-Kentucky Conference"
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
# Row 29
$ws.Range("A29").Value = 59
$ws.Range("C29").Value = "Indianapolis EMS Conference Center"
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 1
# Row 30
$ws.Range("A30").Value = 13
# Row 31
$ws.Range("A31").Value = 10
# Row 33
$ws.Range("A33").Value = 30
# Row 35
$ws.Range("A35").Value = 53
# Row 37
$ws.Range("A37").Value = 36
# Row 40
$ws.Range("A40").Value = 28
# Row 41
$ws.Range("A41").Value = 39
# Row 42
$ws.Range("A42").Value = 32
# Row 47
$ws.Range("A47").Value = 15
$ws.Range("E47").Value = 845
# Row 49
$ws.Range("A49").Value = 9
# Row 51
$ws.Range("A51").Value = 8
# Row 53
$ws.Range("A53").Value = 35
# Row 54
$ws.Range("A54").Value = 40

Write-Host "Applied conference-center location updates"
